$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Capital_Investment
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Capital_Investment")
$ws.Range("B7").Value = 244176.2
$ws.Range("C8").Value = "Rent deposit (7mo @ full rate), Q1 rent (50% rate), service charge, insurance, business rates, legal"

# ---------------------------------------------------------------
# Sheet: Costs_Tracker
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Costs_Tracker")
$ws.Range("A2").Value = "Warehouse - Rent Deposit (7 months @ full rate)"
$ws.Range("B2").Value = 179743.2
$ws.Range("D2").Value = 179743.2
$ws.Range("E2").Value = "£149,786 + VAT (£29,957.20) - Based on full rate £21,398/mo"
$ws.Range("A3").Value = "Warehouse - Q1 Rent (50% reduced rate)"
$ws.Range("E3").Value = "£32,097 + VAT (£6,419.40) - 2026 rate: £10,699/mo"

# ---------------------------------------------------------------
# Sheet: Monthly_Cashflow
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Monthly_Cashflow")

$ws.Range("B4").Value = "Rent Deposit (7mo @ full rate)"
$ws.Range("C4").Value = 149786
$ws.Range("D4").Value = 29957.2
$ws.Range("E4").Value = 179743.2
$ws.Range("G4").Value = 432140.2
$ws.Range("H4").Value = "7 months @ £21,398/mo (2027 rate)"

$ws.Range("B5").Value = "Q1 Rent (50% reduced)"
$ws.Range("G5").Value = 393623.8
$ws.Range("H5").Value = "2026 rate: £10,699/mo"

$ws.Range("G6").Value = 390623.8
$ws.Range("G7").Value = 385823.8
$ws.Range("G8").Value = 380823.8
$ws.Range("G9").Value = 375823.8
$ws.Range("G10").Value = 370823.8
$ws.Range("G11").Value = 365823.8

$ws.Range("B12").Value = "Q2 Rent (50% reduced)"
$ws.Range("G12").Value = 327307.4
$ws.Range("H12").Value = "2026 rate: £10,699/mo"

$ws.Range("G13").Value = 324307.4
$ws.Range("G14").Value = 319307.4
$ws.Range("G15").Value = 314307.4
$ws.Range("G16").Value = 309307.4

$ws.Range("B17").Value = "Q3 Rent (50% reduced)"
$ws.Range("G17").Value = 270791
$ws.Range("H17").Value = "2026 rate: £10,699/mo"

$ws.Range("G18").Value = 267791
